$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new quarter columns before column D; existing D:K data shifts to F:M
$ws.Range("D:E").Insert()

# Copy number format/style from column F (the old column D, now shifted) onto the
# two new columns D:E for every data row so the new cells match the sheet formatting
$ws.Range("F7:F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F8:F8").Copy()
$ws.Range("D8:E8").PasteSpecial(-4122)
$ws.Range("F9:F9").Copy()
$ws.Range("D9:E9").PasteSpecial(-4122)
$ws.Range("F10:F10").Copy()
$ws.Range("D10:E10").PasteSpecial(-4122)
$ws.Range("F11:F11").Copy()
$ws.Range("D11:E11").PasteSpecial(-4122)
$ws.Range("F12:F12").Copy()
$ws.Range("D12:E12").PasteSpecial(-4122)
$ws.Range("F13:F13").Copy()
$ws.Range("D13:E13").PasteSpecial(-4122)
$ws.Range("F14:F14").Copy()
$ws.Range("D14:E14").PasteSpecial(-4122)
$ws.Range("F15:F15").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)
$ws.Range("F16:F16").Copy()
$ws.Range("D16:E16").PasteSpecial(-4122)
$ws.Range("F17:F17").Copy()
$ws.Range("D17:E17").PasteSpecial(-4122)
$ws.Range("F18:F18").Copy()
$ws.Range("D18:E18").PasteSpecial(-4122)
$ws.Range("F19:F19").Copy()
$ws.Range("D19:E19").PasteSpecial(-4122)
$ws.Range("F20:F20").Copy()
$ws.Range("D20:E20").PasteSpecial(-4122)
$ws.Range("F21:F21").Copy()
$ws.Range("D21:E21").PasteSpecial(-4122)
$ws.Range("F22:F22").Copy()
$ws.Range("D22:E22").PasteSpecial(-4122)
$ws.Range("F23:F23").Copy()
$ws.Range("D23:E23").PasteSpecial(-4122)
$ws.Range("F24:F24").Copy()
$ws.Range("D24:E24").PasteSpecial(-4122)
$ws.Range("F25:F25").Copy()
$ws.Range("D25:E25").PasteSpecial(-4122)
$ws.Range("F26:F26").Copy()
$ws.Range("D26:E26").PasteSpecial(-4122)
$ws.Range("F27:F27").Copy()
$ws.Range("D27:E27").PasteSpecial(-4122)
$ws.Range("F28:F28").Copy()
$ws.Range("D28:E28").PasteSpecial(-4122)
$ws.Range("F29:F29").Copy()
$ws.Range("D29:E29").PasteSpecial(-4122)
$ws.Range("F30:F30").Copy()
$ws.Range("D30:E30").PasteSpecial(-4122)
$ws.Range("F31:F31").Copy()
$ws.Range("D31:E31").PasteSpecial(-4122)
$ws.Range("F32:F32").Copy()
$ws.Range("D32:E32").PasteSpecial(-4122)
$ws.Range("F33:F33").Copy()
$ws.Range("D33:E33").PasteSpecial(-4122)
$ws.Range("F34:F34").Copy()
$ws.Range("D34:E34").PasteSpecial(-4122)
$ws.Range("F35:F35").Copy()
$ws.Range("D35:E35").PasteSpecial(-4122)
$ws.Range("F38:F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F39:F39").Copy()
$ws.Range("D39:E39").PasteSpecial(-4122)
$ws.Range("F40:F40").Copy()
$ws.Range("D40:E40").PasteSpecial(-4122)
$ws.Range("F41:F41").Copy()
$ws.Range("D41:E41").PasteSpecial(-4122)
$ws.Range("F42:F42").Copy()
$ws.Range("D42:E42").PasteSpecial(-4122)
$ws.Range("F43:F43").Copy()
$ws.Range("D43:E43").PasteSpecial(-4122)
$ws.Range("F44:F44").Copy()
$ws.Range("D44:E44").PasteSpecial(-4122)
$ws.Range("F45:F45").Copy()
$ws.Range("D45:E45").PasteSpecial(-4122)
$ws.Range("F46:F46").Copy()
$ws.Range("D46:E46").PasteSpecial(-4122)
$ws.Range("F47:F47").Copy()
$ws.Range("D47:E47").PasteSpecial(-4122)
$ws.Range("F48:F48").Copy()
$ws.Range("D48:E48").PasteSpecial(-4122)
$ws.Range("F49:F49").Copy()
$ws.Range("D49:E49").PasteSpecial(-4122)
$ws.Range("F50:F50").Copy()
$ws.Range("D50:E50").PasteSpecial(-4122)
$ws.Range("F51:F51").Copy()
$ws.Range("D51:E51").PasteSpecial(-4122)
$ws.Range("F52:F52").Copy()
$ws.Range("D52:E52").PasteSpecial(-4122)
$ws.Range("F53:F53").Copy()
$ws.Range("D53:E53").PasteSpecial(-4122)
$ws.Range("F54:F54").Copy()
$ws.Range("D54:E54").PasteSpecial(-4122)
$ws.Range("F55:F55").Copy()
$ws.Range("D55:E55").PasteSpecial(-4122)
$ws.Range("F56:F56").Copy()
$ws.Range("D56:E56").PasteSpecial(-4122)
$ws.Range("F57:F57").Copy()
$ws.Range("D57:E57").PasteSpecial(-4122)
$ws.Range("F58:F58").Copy()
$ws.Range("D58:E58").PasteSpecial(-4122)
$ws.Range("F59:F59").Copy()
$ws.Range("D59:E59").PasteSpecial(-4122)
$ws.Range("F60:F60").Copy()
$ws.Range("D60:E60").PasteSpecial(-4122)
$ws.Range("F61:F61").Copy()
$ws.Range("D61:E61").PasteSpecial(-4122)
$ws.Range("F62:F62").Copy()
$ws.Range("D62:E62").PasteSpecial(-4122)
$ws.Range("F63:F63").Copy()
$ws.Range("D63:E63").PasteSpecial(-4122)
$ws.Range("F64:F64").Copy()
$ws.Range("D64:E64").PasteSpecial(-4122)
$ws.Range("F65:F65").Copy()
$ws.Range("D65:E65").PasteSpecial(-4122)
$ws.Range("F66:F66").Copy()
$ws.Range("D66:E66").PasteSpecial(-4122)
$ws.Range("F67:F67").Copy()
$ws.Range("D67:E67").PasteSpecial(-4122)
$ws.Range("F68:F68").Copy()
$ws.Range("D68:E68").PasteSpecial(-4122)
$ws.Range("F69:F69").Copy()
$ws.Range("D69:E69").PasteSpecial(-4122)
$ws.Range("F70:F70").Copy()
$ws.Range("D70:E70").PasteSpecial(-4122)
$ws.Range("F71:F71").Copy()
$ws.Range("D71:E71").PasteSpecial(-4122)
$ws.Range("F72:F72").Copy()
$ws.Range("D72:E72").PasteSpecial(-4122)
$ws.Range("F73:F73").Copy()
$ws.Range("D73:E73").PasteSpecial(-4122)
$ws.Range("F74:F74").Copy()
$ws.Range("D74:E74").PasteSpecial(-4122)
$ws.Range("F75:F75").Copy()
$ws.Range("D75:E75").PasteSpecial(-4122)
$ws.Range("F76:F76").Copy()
$ws.Range("D76:E76").PasteSpecial(-4122)
$ws.Range("F77:F77").Copy()
$ws.Range("D77:E77").PasteSpecial(-4122)
$ws.Range("F80:F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$ws.Range("F81:F81").Copy()
$ws.Range("D81:E81").PasteSpecial(-4122)
$ws.Range("F82:F82").Copy()
$ws.Range("D82:E82").PasteSpecial(-4122)
$ws.Range("F83:F83").Copy()
$ws.Range("D83:E83").PasteSpecial(-4122)
$ws.Range("F84:F84").Copy()
$ws.Range("D84:E84").PasteSpecial(-4122)
$ws.Range("F85:F85").Copy()
$ws.Range("D85:E85").PasteSpecial(-4122)
$ws.Range("F86:F86").Copy()
$ws.Range("D86:E86").PasteSpecial(-4122)
$ws.Range("F87:F87").Copy()
$ws.Range("D87:E87").PasteSpecial(-4122)
$ws.Range("F88:F88").Copy()
$ws.Range("D88:E88").PasteSpecial(-4122)
$ws.Range("F89:F89").Copy()
$ws.Range("D89:E89").PasteSpecial(-4122)
$ws.Range("F90:F90").Copy()
$ws.Range("D90:E90").PasteSpecial(-4122)
$ws.Range("F91:F91").Copy()
$ws.Range("D91:E91").PasteSpecial(-4122)
$ws.Range("F92:F92").Copy()
$ws.Range("D92:E92").PasteSpecial(-4122)
$ws.Range("F93:F93").Copy()
$ws.Range("D93:E93").PasteSpecial(-4122)
$ws.Range("F94:F94").Copy()
$ws.Range("D94:E94").PasteSpecial(-4122)
$ws.Range("F95:F95").Copy()
$ws.Range("D95:E95").PasteSpecial(-4122)
$ws.Range("F96:F96").Copy()
$ws.Range("D96:E96").PasteSpecial(-4122)
$ws.Range("F97:F97").Copy()
$ws.Range("D97:E97").PasteSpecial(-4122)
$ws.Range("F98:F98").Copy()
$ws.Range("D98:E98").PasteSpecial(-4122)
$ws.Range("F99:F99").Copy()
$ws.Range("D99:E99").PasteSpecial(-4122)
$ws.Range("F100:F100").Copy()
$ws.Range("D100:E100").PasteSpecial(-4122)
$ws.Range("F101:F101").Copy()
$ws.Range("D101:E101").PasteSpecial(-4122)
$ws.Range("F102:F102").Copy()
$ws.Range("D102:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate final values for every quarter column (D:M) on every data row, reflecting
# both the two newly added quarters and the restated figures for prior quarters
# Row 7
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
# Row 8
$ws.Range("D8").Value = 520500
$ws.Range("E8").Value = 556100
$ws.Range("F8").Value = 1094000
$ws.Range("G8").Value = 732400
$ws.Range("H8").Value = 518700
$ws.Range("I8").Value = 496600
$ws.Range("J8").Value = 1123600
$ws.Range("K8").Value = 756300
$ws.Range("L8").Value = 659300
$ws.Range("M8").Value = 738200
# Row 9
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "NA"
$ws.Range("H9").Value = "NA"
$ws.Range("I9").Value = "NA"
$ws.Range("J9").Value = "NA"
$ws.Range("K9").Value = "NA"
$ws.Range("L9").Value = "NA"
$ws.Range("M9").Value = "NA"
# Row 10
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "NA"
$ws.Range("H10").Value = "NA"
$ws.Range("I10").Value = "NA"
$ws.Range("J10").Value = "NA"
$ws.Range("K10").Value = "NA"
$ws.Range("L10").Value = "NA"
$ws.Range("M10").Value = "NA"
# Row 11
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
# Row 12
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
# Row 14
$ws.Range("D14").Value = 5400
$ws.Range("E14").Value = 500
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = "NA"
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
# Row 15
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = 500
$ws.Range("G15").Value = 500
$ws.Range("H15").Value = 500
$ws.Range("I15").Value = 500
$ws.Range("J15").Value = 1100
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 600
$ws.Range("M15").Value = 600
# Row 16
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
# Row 17
$ws.Range("D17").Value = 736300
$ws.Range("E17").Value = 787300
$ws.Range("F17").Value = 1082000
$ws.Range("G17").Value = 702400
$ws.Range("H17").Value = 635700
$ws.Range("I17").Value = 531600
$ws.Range("J17").Value = 1084900
$ws.Range("K17").Value = 720500
$ws.Range("L17").Value = 726600
$ws.Range("M17").Value = 691000
# Row 18
$ws.Range("D18").Value = -215800
$ws.Range("E18").Value = -231200
$ws.Range("F18").Value = 12000
$ws.Range("G18").Value = 30100
$ws.Range("H18").Value = -117000
$ws.Range("I18").Value = -35000
$ws.Range("J18").Value = 38700
$ws.Range("K18").Value = 35800
$ws.Range("L18").Value = -67300
$ws.Range("M18").Value = 47200
# Row 19
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = ""
$ws.Range("J19").Value = ""
$ws.Range("K19").Value = ""
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = ""
# Row 20
$ws.Range("D20").Value = 2600
$ws.Range("E20").Value = -500
$ws.Range("F20").Value = 2400
$ws.Range("G20").Value = -2400
$ws.Range("H20").Value = -2800
$ws.Range("I20").Value = -3600
$ws.Range("J20").Value = -11500
$ws.Range("K20").Value = -1900
$ws.Range("L20").Value = 5100
$ws.Range("M20").Value = 700
# Row 21
$ws.Range("D21").Value = -211800
$ws.Range("E21").Value = -232600
$ws.Range("F21").Value = 20000
$ws.Range("G21").Value = 31200
$ws.Range("H21").Value = -119800
$ws.Range("I21").Value = -37600
$ws.Range("J21").Value = 34700
$ws.Range("K21").Value = 35500
$ws.Range("L21").Value = -58700
$ws.Range("M21").Value = 52900
# Row 22
$ws.Range("D22").Value = 4800
$ws.Range("E22").Value = 4800
$ws.Range("F22").Value = 9700
$ws.Range("G22").Value = 4800
$ws.Range("H22").Value = 4800
$ws.Range("I22").Value = 4800
$ws.Range("J22").Value = 13600
$ws.Range("K22").Value = 6900
$ws.Range("L22").Value = 6900
$ws.Range("M22").Value = 6900
# Row 23
$ws.Range("D23").Value = -218000
$ws.Range("E23").Value = -236600
$ws.Range("F23").Value = 4800
$ws.Range("G23").Value = 22800
$ws.Range("H23").Value = -124600
$ws.Range("I23").Value = -43400
$ws.Range("J23").Value = 13600
$ws.Range("K23").Value = 27000
$ws.Range("L23").Value = -69000
$ws.Range("M23").Value = 41000
# Row 24
$ws.Range("D24").Value = 1400
$ws.Range("E24").Value = -7400
$ws.Range("F24").Value = 6500
$ws.Range("G24").Value = 500
$ws.Range("H24").Value = -4800
$ws.Range("I24").Value = 1700
$ws.Range("J24").Value = -3700
$ws.Range("K24").Value = 500
$ws.Range("L24").Value = 400
$ws.Range("M24").Value = 200
# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
# Row 26
$ws.Range("D26").Value = -219400
$ws.Range("E26").Value = -229100
$ws.Range("F26").Value = -1700
$ws.Range("G26").Value = 22300
$ws.Range("H26").Value = -119800
$ws.Range("I26").Value = -45100
$ws.Range("J26").Value = 17300
$ws.Range("K26").Value = 26500
$ws.Range("L26").Value = -69400
$ws.Range("M26").Value = 40800
# Row 27
$ws.Range("D27").Value = -219500
$ws.Range("E27").Value = -237700
$ws.Range("F27").Value = -19000
$ws.Range("G27").Value = 13700
$ws.Range("H27").Value = -128600
$ws.Range("I27").Value = -53700
$ws.Range("J27").Value = 5300
$ws.Range("K27").Value = 20500
$ws.Range("L27").Value = -74700
$ws.Range("M27").Value = 31800
# Row 28
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
# Row 29
$ws.Range("D29").Value = -49800
$ws.Range("E29").Value = -71100
$ws.Range("F29").Value = 26800
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = -100
$ws.Range("I29").Value = -9900
$ws.Range("J29").Value = -7200
$ws.Range("K29").Value = "NA"
$ws.Range("L29").Value = "NA"
$ws.Range("M29").Value = "NA"
# Row 30
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
# Row 31
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
# Row 32
$ws.Range("D32").Value = -2600
$ws.Range("E32").Value = 500
$ws.Range("F32").Value = -2400
$ws.Range("G32").Value = 2400
$ws.Range("H32").Value = 2800
$ws.Range("I32").Value = 3600
$ws.Range("J32").Value = 11500
$ws.Range("K32").Value = 1900
$ws.Range("L32").Value = -5100
$ws.Range("M32").Value = -700
# Row 33
$ws.Range("D33").Value = -269200
$ws.Range("E33").Value = -308800
$ws.Range("F33").Value = 7800
$ws.Range("G33").Value = 13700
$ws.Range("H33").Value = -128600
$ws.Range("I33").Value = -63600
$ws.Range("J33").Value = -1900
$ws.Range("K33").Value = 20500
$ws.Range("L33").Value = -74700
$ws.Range("M33").Value = 31800
# Row 34
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
# Row 35
$ws.Range("D35").Value = -269200
$ws.Range("E35").Value = -308800
$ws.Range("F35").Value = 7800
$ws.Range("G35").Value = 13700
$ws.Range("H35").Value = -128600
$ws.Range("I35").Value = -63600
$ws.Range("J35").Value = -1900
$ws.Range("K35").Value = 20500
$ws.Range("L35").Value = -74700
$ws.Range("M35").Value = 31800
# Row 38
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
# Row 39
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = ""
$ws.Range("F39").Value = ""
$ws.Range("G39").Value = ""
$ws.Range("H39").Value = ""
$ws.Range("I39").Value = ""
$ws.Range("J39").Value = ""
$ws.Range("K39").Value = ""
$ws.Range("L39").Value = ""
$ws.Range("M39").Value = ""
# Row 40
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = ""
$ws.Range("F40").Value = ""
$ws.Range("G40").Value = ""
$ws.Range("H40").Value = ""
$ws.Range("I40").Value = ""
$ws.Range("J40").Value = ""
$ws.Range("K40").Value = ""
$ws.Range("L40").Value = ""
$ws.Range("M40").Value = ""
# Row 41
$ws.Range("D41").Value = 200800
$ws.Range("E41").Value = 94600
$ws.Range("F41").Value = 72100
$ws.Range("G41").Value = 73300
$ws.Range("H41").Value = 54500
$ws.Range("I41").Value = 182700
$ws.Range("J41").Value = 246800
$ws.Range("K41").Value = 83500
$ws.Range("L41").Value = 45700
$ws.Range("M41").Value = 67500
# Row 42
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
# Row 43
$ws.Range("D43").Value = 67300
$ws.Range("E43").Value = 161400
$ws.Range("F43").Value = 467600
$ws.Range("G43").Value = 548200
$ws.Range("H43").Value = 72500
$ws.Range("I43").Value = 479500
$ws.Range("J43").Value = 499800
$ws.Range("K43").Value = 615600
$ws.Range("L43").Value = 410200
$ws.Range("M43").Value = 535000
# Row 44
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
# Row 45
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
# Row 46
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
# Row 47
$ws.Range("D47").Value = 4258900
$ws.Range("E47").Value = 4005000
$ws.Range("F47").Value = 5179900
$ws.Range("G47").Value = 5230500
$ws.Range("H47").Value = 3979900
$ws.Range("I47").Value = 5178000
$ws.Range("J47").Value = 5021200
$ws.Range("K47").Value = 4900000
$ws.Range("L47").Value = 4904900
$ws.Range("M47").Value = 4688600
# Row 48
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
# Row 49
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("F49").Value = 74700
$ws.Range("G49").Value = 75100
$ws.Range("H49").Value = 75600
$ws.Range("I49").Value = 76100
$ws.Range("J49").Value = 76600
$ws.Range("K49").Value = 77200
$ws.Range("L49").Value = 77700
$ws.Range("M49").Value = 80100
# Row 50
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
# Row 51
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
# Row 52
$ws.Range("D52").Value = 304600
$ws.Range("E52").Value = 1785500
$ws.Range("F52").Value = 154700
$ws.Range("G52").Value = 87800
$ws.Range("H52").Value = 1996700
$ws.Range("I52").Value = 131600
$ws.Range("J52").Value = 190500
$ws.Range("K52").Value = 109700
$ws.Range("L52").Value = 103800
$ws.Range("M52").Value = 366600
# Row 53
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
# Row 54
$ws.Range("D54").Value = 5287500
$ws.Range("E54").Value = 6536500
$ws.Range("F54").Value = 6718200
$ws.Range("G54").Value = 6757300
$ws.Range("H54").Value = 6644200
$ws.Range("I54").Value = 6839100
$ws.Range("J54").Value = 6791200
$ws.Range("K54").Value = 6557000
$ws.Range("L54").Value = 6252300
$ws.Range("M54").Value = 6468000
# Row 55
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = ""
$ws.Range("F55").Value = ""
$ws.Range("G55").Value = ""
$ws.Range("H55").Value = ""
$ws.Range("I55").Value = ""
$ws.Range("J55").Value = ""
$ws.Range("K55").Value = ""
$ws.Range("L55").Value = ""
$ws.Range("M55").Value = ""
# Row 56
$ws.Range("D56").Value = ""
$ws.Range("E56").Value = ""
$ws.Range("F56").Value = ""
$ws.Range("G56").Value = ""
$ws.Range("H56").Value = ""
$ws.Range("I56").Value = ""
$ws.Range("J56").Value = ""
$ws.Range("K56").Value = ""
$ws.Range("L56").Value = ""
$ws.Range("M56").Value = ""
# Row 57
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 0
# Row 58
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
# Row 59
$ws.Range("D59").Value = 4321900
$ws.Range("E59").Value = 4169100
$ws.Range("F59").Value = 5357300
$ws.Range("G59").Value = 5344300
$ws.Range("H59").Value = 3707700
$ws.Range("I59").Value = 5163300
$ws.Range("J59").Value = 5036100
$ws.Range("K59").Value = 4835400
$ws.Range("L59").Value = 4539700
$ws.Range("M59").Value = 4560300
# Row 60
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 0
# Row 61
$ws.Range("D61").Value = 254700
$ws.Range("E61").Value = 254600
$ws.Range("F61").Value = 254600
$ws.Range("G61").Value = 254500
$ws.Range("H61").Value = 254500
$ws.Range("I61").Value = 254400
$ws.Range("J61").Value = 254400
$ws.Range("K61").Value = 351500
$ws.Range("L61").Value = 351400
$ws.Range("M61").Value = 351300
# Row 62
$ws.Range("D62").Value = 156000
$ws.Range("E62").Value = 1339600
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "NA"
$ws.Range("H62").Value = 1449400
$ws.Range("I62").Value = "NA"
$ws.Range("J62").Value = "NA"
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 0
# Row 63
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
# Row 64
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
# Row 65
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
# Row 66
$ws.Range("D66").Value = 4733200
$ws.Range("E66").Value = 5763900
$ws.Range("F66").Value = 5612400
$ws.Range("G66").Value = 5599400
$ws.Range("H66").Value = 5412000
$ws.Range("I66").Value = 5418100
$ws.Range("J66").Value = 5290800
$ws.Range("K66").Value = 5187200
$ws.Range("L66").Value = 4891500
$ws.Range("M66").Value = 4912700
# Row 67
$ws.Range("D67").Value = ""
$ws.Range("E67").Value = ""
$ws.Range("F67").Value = ""
$ws.Range("G67").Value = ""
$ws.Range("H67").Value = ""
$ws.Range("I67").Value = ""
$ws.Range("J67").Value = ""
$ws.Range("K67").Value = ""
$ws.Range("L67").Value = ""
$ws.Range("M67").Value = ""
# Row 68
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
# Row 69
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
# Row 70
$ws.Range("D70").Value = 465000
$ws.Range("E70").Value = 465000
$ws.Range("F70").Value = 465000
$ws.Range("G70").Value = 465000
$ws.Range("H70").Value = 465000
$ws.Range("I70").Value = 465000
$ws.Range("J70").Value = 465000
$ws.Range("K70").Value = 315000
$ws.Range("L70").Value = 315000
$ws.Range("M70").Value = 315000
# Row 71
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
# Row 72
$ws.Range("D72").Value = -563900
$ws.Range("E72").Value = -294700
$ws.Range("F72").Value = 18300
$ws.Range("G72").Value = 36700
$ws.Range("H72").Value = 35500
$ws.Range("I72").Value = 181500
$ws.Range("J72").Value = 257800
$ws.Range("K72").Value = 293200
$ws.Range("L72").Value = 285700
$ws.Range("M72").Value = 373300
# Row 73
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
# Row 74
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
# Row 75
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
# Row 76
$ws.Range("D76").Value = 89300
$ws.Range("E76").Value = 307600
$ws.Range("F76").Value = 640700
$ws.Range("G76").Value = 692900
$ws.Range("H76").Value = 767200
$ws.Range("I76").Value = 956000
$ws.Range("J76").Value = 1035400
$ws.Range("K76").Value = 1054800
$ws.Range("L76").Value = 1045800
$ws.Range("M76").Value = 1240300
# Row 77
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
# Row 80
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
# Row 81
$ws.Range("D81").Value = -269200
$ws.Range("E81").Value = -308800
$ws.Range("F81").Value = 7800
$ws.Range("G81").Value = 13700
$ws.Range("H81").Value = -128600
$ws.Range("I81").Value = -63600
$ws.Range("J81").Value = -1900
$ws.Range("K81").Value = 20500
$ws.Range("L81").Value = -74700
$ws.Range("M81").Value = 31800
# Row 82
$ws.Range("D82").Value = ""
$ws.Range("E82").Value = ""
$ws.Range("F82").Value = ""
$ws.Range("G82").Value = ""
$ws.Range("H82").Value = ""
$ws.Range("I82").Value = ""
$ws.Range("J82").Value = ""
$ws.Range("K82").Value = ""
$ws.Range("L82").Value = ""
$ws.Range("M82").Value = ""
# Row 83
$ws.Range("D83").Value = 1400
$ws.Range("E83").Value = -800
$ws.Range("F83").Value = 5600
$ws.Range("G83").Value = 3500
$ws.Range("H83").Value = 1400
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 7500
$ws.Range("K83").Value = 1600
$ws.Range("L83").Value = 3500
$ws.Range("M83").Value = 5000
# Row 84
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
# Row 85
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
# Row 86
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
# Row 87
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
# Row 88
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
# Row 89
$ws.Range("D89").Value = 13600
$ws.Range("E89").Value = 106500
$ws.Range("F89").Value = 56100
$ws.Range("G89").Value = -32100
$ws.Range("H89").Value = 83700
$ws.Range("I89").Value = 151300
$ws.Range("J89").Value = 223500
$ws.Range("K89").Value = 45500
$ws.Range("L89").Value = 144000
$ws.Range("M89").Value = 161100
# Row 90
$ws.Range("D90").Value = ""
$ws.Range("E90").Value = ""
$ws.Range("F90").Value = ""
$ws.Range("G90").Value = ""
$ws.Range("H90").Value = ""
$ws.Range("I90").Value = ""
$ws.Range("J90").Value = ""
$ws.Range("K90").Value = ""
$ws.Range("L90").Value = ""
$ws.Range("M90").Value = ""
# Row 91
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 0
# Row 92
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
# Row 93
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
# Row 94
$ws.Range("D94").Value = 5700
$ws.Range("E94").Value = 5000
$ws.Range("F94").Value = 22200
$ws.Range("G94").Value = 22000
$ws.Range("H94").Value = -174800
$ws.Range("I94").Value = -239300
$ws.Range("J94").Value = 54300
$ws.Range("K94").Value = 16800
$ws.Range("L94").Value = -147800
$ws.Range("M94").Value = -202000
# Row 95
$ws.Range("D95").Value = ""
$ws.Range("E95").Value = ""
$ws.Range("F95").Value = ""
$ws.Range("G95").Value = ""
$ws.Range("H95").Value = ""
$ws.Range("I95").Value = ""
$ws.Range("J95").Value = ""
$ws.Range("K95").Value = ""
$ws.Range("L95").Value = ""
$ws.Range("M95").Value = ""
# Row 96
$ws.Range("D96").Value = -4200
$ws.Range("E96").Value = -12500
$ws.Range("F96").Value = -24900
$ws.Range("G96").Value = -12500
$ws.Range("H96").Value = -12700
$ws.Range("I96").Value = -13000
$ws.Range("J96").Value = -25900
$ws.Range("K96").Value = -12900
$ws.Range("L96").Value = -12100
$ws.Range("M96").Value = -10400
# Row 97
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
# Row 98
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
# Row 99
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
# Row 100
$ws.Range("D100").Value = -4100
$ws.Range("E100").Value = -21700
$ws.Range("F100").Value = -42200
$ws.Range("G100").Value = -21200
$ws.Range("H100").Value = -32000
$ws.Range("I100").Value = -36000
$ws.Range("J100").Value = 7600
$ws.Range("K100").Value = -19000
$ws.Range("L100").Value = -17100
$ws.Range("M100").Value = -19300
# Row 101
$ws.Range("D101").Value = -400
$ws.Range("E101").Value = -300
$ws.Range("F101").Value = -800
$ws.Range("G101").Value = 800
$ws.Range("H101").Value = 300
$ws.Range("I101").Value = 1000
$ws.Range("J101").Value = 2400
$ws.Range("K101").Value = 500
$ws.Range("L101").Value = -800
$ws.Range("M101").Value = 200
# Row 102
$ws.Range("D102").Value = 14700
$ws.Range("E102").Value = 89500
$ws.Range("F102").Value = 35300
$ws.Range("G102").Value = -30500
$ws.Range("H102").Value = -142600
$ws.Range("I102").Value = -123000
$ws.Range("J102").Value = 287800
$ws.Range("K102").Value = 43700
$ws.Range("L102").Value = -21700
$ws.Range("M102").Value = -60000
